$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: copy a data row (columns A..I) from $srcRow to $dstRow while
# preserving each cell's original type (text vs number) by using
# Range.Copy for non-empty cells (rather than re-typing scalar .Value2,
# which would coerce numeric-looking text like "4" into a real number)
# and ClearContents for cells that must end up empty (rather than
# Range.Clear, which leaves a stray empty <c> element behind).
function Copy-DataRow($srcRow, $dstRow) {
    for ($c = 1; $c -le 9; $c++) {
        $srcCell = $ws.Cells.Item($srcRow, $c)
        $dstCell = $ws.Cells.Item($dstRow, $c)
        $v = $srcCell.Value2
        if ($v -eq $null -or $v -eq "") {
            $dstCell.ClearContents()
        } else {
            $srcCell.Copy($dstCell)
        }
    }
}

# Swap the content of rows 5 and 6 using row 100 as scratch space.
Copy-DataRow 5 100
Copy-DataRow 6 5
Copy-DataRow 100 6
$ws.Range("A100:I100").ClearContents()

# Add a space before the colon in the statut_name values throughout
# column B, e.g. "4: pas de ..." -> "4 : pas de ..." and
# "2: résultats ..." -> "2 : résultats ...".
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace '^(\d+):', '$1 :'
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
